$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Mark "OK" in column G for a few rows that were missing it (row 3, 4, 5, 27)
$ws.Range("G3").Value = "OK"
$ws.Range("G4").Value = "OK"
$ws.Range("G5").Value = "OK"
$ws.Range("G27").Value = "OK"

# Append two new test rows (29, 30 -> sheet rows 30, 31)
$ws.Range("A30").Value = 29
$ws.Range("B30").Value = "передача пустого аргумента в функцию, выбор его в селекте"

$ws.Range("A31").Value = 30
$ws.Range("B31").Value = "передача некорректного аргумента в каунтер"

# Extra blank placeholder rows (31 .. 43) in column A only, continuing the numbering
$ws.Range("A32").Value = 31
$ws.Range("A33").Value = 32
$ws.Range("A34").Value = 33
$ws.Range("A35").Value = 34
$ws.Range("A36").Value = 35
$ws.Range("A37").Value = 36
$ws.Range("A38").Value = 37
$ws.Range("A39").Value = 38
$ws.Range("A40").Value = 39
$ws.Range("A41").Value = 40
$ws.Range("A42").Value = 41
$ws.Range("A43").Value = 42
$ws.Range("A44").Value = 43

# Scroll the view down/right a bit and move the active selection, like a user
# who just finished adding the new rows further down the sheet.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 8
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B32").Select()
